$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows1 = @(4,5,6,7,8,9,10,12,13,14,15,21,22)
foreach ($r in $rows1) {
    for ($c = 10; $c -le 45; $c++) {
        $ws.Cells.Item($r, $c).Value = 1
    }
}

# Row 96
$ws.Cells.Item(96, 10).Value = 17654249.7138674
$ws.Cells.Item(96, 11).Value = 14052415.3315779
$ws.Cells.Item(96, 12).Value = 13820024.81106
$ws.Cells.Item(96, 13).Value = 18419922.1048115
$ws.Cells.Item(96, 14).Value = 15265204.5558733
$ws.Cells.Item(96, 15).Value = 12093689.4248959
$ws.Cells.Item(96, 16).Value = 12093690
$ws.Cells.Item(96, 17).Value = 12982635.5708305
$ws.Cells.Item(96, 18).Value = 13936923.0040619
$ws.Cells.Item(96, 19).Value = 14961355.2472786
$ws.Cells.Item(96, 20).Value = 16061088.2882851
$ws.Cells.Item(96, 21).Value = 17241657.1052954
$ws.Cells.Item(96, 22).Value = 18509003.5245877
$ws.Cells.Item(96, 23).Value = 19869506.1258342
$ws.Cells.Item(96, 24).Value = 21330012.3456191
$ws.Cells.Item(96, 25).Value = 22897872.9407227
$ws.Cells.Item(96, 26).Value = 24580978.9846262
$ws.Cells.Item(96, 27).Value = 26387801.5834411
$ws.Cells.Item(96, 28).Value = 28327434.5111542
$ws.Cells.Item(96, 29).Value = 30409639.9787726
$ws.Cells.Item(96, 30).Value = 32644897.7677254
$ws.Cells.Item(96, 31).Value = 35044457.9748114
$ws.Cells.Item(96, 32).Value = 37620397.6341598
$ws.Cells.Item(96, 33).Value = 40385681.5011822
$ws.Cells.Item(96, 34).Value = 43354227.3044439
$ws.Cells.Item(96, 35).Value = 46540975.7938681
$ws.Cells.Item(96, 36).Value = 49961965.9378266
$ws.Cells.Item(96, 37).Value = 53634415.6475857
$ws.Cells.Item(96, 38).Value = 57576808.4353952
$ws.Cells.Item(96, 39).Value = 61808986.4423726
$ws.Cells.Item(96, 40).Value = 66352250.3043925
$ws.Cells.Item(96, 41).Value = 71229466.3586098
$ws.Cells.Item(96, 42).Value = 76465181.7301885
$ws.Cells.Item(96, 43).Value = 82085747.8784693
$ws.Cells.Item(96, 44).Value = 88119453.2243874
$ws.Cells.Item(96, 45).Value = 94596665.5266563

# Row 97
$ws.Cells.Item(97, 10).Value = 852872.099516197
$ws.Cells.Item(97, 11).Value = 860298.351874753
$ws.Cells.Item(97, 12).Value = 929757.423772548
$ws.Cells.Item(97, 13).Value = 915785.743949109
$ws.Cells.Item(97, 14).Value = 950156.14250452
$ws.Cells.Item(97, 15).Value = 939774.902194865
$ws.Cells.Item(97, 16).Value = 939774.9
$ws.Cells.Item(97, 17).Value = 979711.433935985
$ws.Cells.Item(97, 18).Value = 1021345.10485958
$ws.Cells.Item(97, 19).Value = 1064748.03405101
$ws.Cells.Item(97, 20).Value = 1109995.40764563
$ws.Cells.Item(97, 21).Value = 1157165.60687762
$ws.Cells.Item(97, 22).Value = 1206340.34385839
$ws.Cells.Item(97, 23).Value = 1257604.80312501
$ws.Cells.Item(97, 24).Value = 1311047.78920398
$ws.Cells.Item(97, 25).Value = 1366761.88044566
$ws.Cells.Item(97, 26).Value = 1424843.58939621
$ws.Cells.Item(97, 27).Value = 1485393.52998452
$ws.Cells.Item(97, 28).Value = 1548516.59181403
$ws.Cells.Item(97, 29).Value = 1614322.12186108
$ws.Cells.Item(97, 30).Value = 1682924.11389483
$ws.Cells.Item(97, 31).Value = 1754441.40594663
$ws.Cells.Item(97, 32).Value = 1828997.88617108
$ws.Cells.Item(97, 33).Value = 1906722.70745532
$ws.Cells.Item(97, 34).Value = 1987750.51114831
$ws.Cells.Item(97, 35).Value = 2072221.6602977
$ws.Cells.Item(97, 36).Value = 2160282.48279824
$ws.Cells.Item(97, 37).Value = 2252085.52487308
$ws.Cells.Item(97, 38).Value = 2347789.81532693
$ws.Cells.Item(97, 39).Value = 2447561.14102882
$ws.Cells.Item(97, 40).Value = 2551572.33410186
$ws.Cells.Item(97, 41).Value = 2660003.57131726
$ws.Cells.Item(97, 42).Value = 2773042.68621143
$ws.Cells.Item(97, 43).Value = 2890885.4944667
$ws.Cells.Item(97, 44).Value = 3013736.13311944
$ws.Cells.Item(97, 45).Value = 3141807.41418305

# Row 98
$ws.Cells.Item(98, 10).Value = 491353.270110164
$ws.Cells.Item(98, 11).Value = 453230.009368236
$ws.Cells.Item(98, 12).Value = 469199.451919957
$ws.Cells.Item(98, 13).Value = 441982.475839454
$ws.Cells.Item(98, 14).Value = 457289.149611657
$ws.Cells.Item(98, 15).Value = 429809.496763082
$ws.Cells.Item(98, 16).Value = 429809.5
$ws.Cells.Item(98, 17).Value = 444278.142383931
$ws.Cells.Item(98, 18).Value = 459233.841504471
$ws.Cells.Item(98, 19).Value = 474692.993113095
$ws.Cells.Item(98, 20).Value = 490672.544890129
$ws.Cells.Item(98, 21).Value = 507190.015024289
$ws.Cells.Item(98, 22).Value = 524263.51141766
$ws.Cells.Item(98, 23).Value = 541911.751537169
$ws.Cells.Item(98, 24).Value = 560154.082934312
$ws.Cells.Item(98, 25).Value = 579010.504455649
$ws.Cells.Item(98, 26).Value = 598501.688167289
$ws.Cells.Item(98, 27).Value = 618649.002017428
$ws.Cells.Item(98, 28).Value = 639474.533261771
$ws.Cells.Item(98, 29).Value = 661001.112677525
$ws.Cells.Item(98, 30).Value = 683252.339592498
$ws.Cells.Item(98, 31).Value = 706252.607756761
$ws.Cells.Item(98, 32).Value = 730027.132085216
$ws.Cells.Item(98, 33).Value = 754601.976300403
$ws.Cells.Item(98, 34).Value = 780004.081505843
$ws.Cells.Item(98, 35).Value = 806261.295721243
$ws.Cells.Item(98, 36).Value = 833402.404411943
$ws.Cells.Item(98, 37).Value = 861457.162046067
$ws.Cells.Item(98, 38).Value = 890456.324713994
$ws.Cells.Item(98, 39).Value = 920431.683845879
$ws.Cells.Item(98, 40).Value = 951416.101064216
$ws.Cells.Item(98, 41).Value = 983443.544209636
$ws.Cells.Item(98, 42).Value = 1016549.12457944
$ws.Cells.Item(98, 43).Value = 1050769.13541968
$ws.Cells.Item(98, 44).Value = 1086141.09171301
$ws.Cells.Item(98, 45).Value = 1122703.77130592

# Row 99
$ws.Cells.Item(99, 10).Value = 72764.1102957194
$ws.Cells.Item(99, 11).Value = 70201.0658287004
$ws.Cells.Item(99, 12).Value = 71386.4490428062
$ws.Cells.Item(99, 13).Value = 79583.3540969824
$ws.Cells.Item(99, 14).Value = 78440.0908814858
$ws.Cells.Item(99, 15).Value = 66995.5216394048
$ws.Cells.Item(99, 16).Value = 66995.52
$ws.Cells.Item(99, 17).Value = 69214.9073624498
$ws.Cells.Item(99, 18).Value = 71507.8172569227
$ws.Cells.Item(99, 19).Value = 73876.6852937165
$ws.Cells.Item(99, 20).Value = 76324.0277685649
$ws.Cells.Item(99, 21).Value = 78852.4443355357
$ws.Cells.Item(99, 22).Value = 81464.6207684757
$ws.Cells.Item(99, 23).Value = 84163.3318139355
$ws.Cells.Item(99, 24).Value = 86951.4441386031
$ws.Cells.Item(99, 25).Value = 89831.9193743797
$ws.Cells.Item(99, 26).Value = 92807.8172643298
$ws.Cells.Item(99, 27).Value = 95882.298912849
$ws.Cells.Item(99, 28).Value = 99058.6301435017
$ws.Cells.Item(99, 29).Value = 102340.184968094
$ws.Cells.Item(99, 30).Value = 105730.44917067
$ws.Cells.Item(99, 31).Value = 109233.024010234
$ws.Cells.Item(99, 32).Value = 112851.630046138
$ws.Cells.Item(99, 33).Value = 116590.111090188
$ws.Cells.Item(99, 34).Value = 120452.438289681
$ws.Cells.Item(99, 35).Value = 124442.714345697
$ws.Cells.Item(99, 36).Value = 128565.177871136
$ws.Cells.Item(99, 37).Value = 132824.207893118
$ws.Cells.Item(99, 38).Value = 137224.328504547
$ws.Cells.Item(99, 39).Value = 141770.213669759
$ws.Cells.Item(99, 40).Value = 146466.692189374
$ws.Cells.Item(99, 41).Value = 151318.752829622
$ws.Cells.Item(99, 42).Value = 156331.549621582
$ws.Cells.Item(99, 43).Value = 161510.407335984
$ws.Cells.Item(99, 44).Value = 166860.827139362
$ws.Cells.Item(99, 45).Value = 172388.492437594

# Row 100
$ws.Cells.Item(100, 10).Value = 265449.737103148
$ws.Cells.Item(100, 11).Value = 133762.956322085
$ws.Cells.Item(100, 12).Value = 342546.377465837
$ws.Cells.Item(100, 13).Value = 321850.093104366
$ws.Cells.Item(100, 14).Value = 304542.555223962
$ws.Cells.Item(100, 15).Value = 329257.541195101
$ws.Cells.Item(100, 16).Value = 329257.5
$ws.Cells.Item(100, 17).Value = 587980.17819684
$ws.Cells.Item(100, 18).Value = 1050000.95655342
$ws.Cells.Item(100, 19).Value = 1875066.62579025
$ws.Cells.Item(100, 20).Value = 3348449.18874467
$ws.Cells.Item(100, 21).Value = 5979580.57350601
$ws.Cells.Item(100, 22).Value = 10678192.1479463
$ws.Cells.Item(100, 23).Value = 19068860.4571484
$ws.Cells.Item(100, 24).Value = 34052715.4874371
$ws.Cells.Item(100, 25).Value = 60810525.8661976
$ws.Cells.Item(100, 26).Value = 108593984.450014
$ws.Cells.Item(100, 27).Value = 193924543.337734
$ws.Cells.Item(100, 28).Value = 346305817.023034
$ws.Cells.Item(100, 29).Value = 618424655.486377
$ws.Cells.Item(100, 30).Value = 1104367976.8394
$ws.Cells.Item(100, 31).Value = 1972153951.89757
$ws.Cells.Item(100, 32).Value = 3521825416.48508
$ws.Cells.Item(100, 33).Value = 6289191699.39352
$ws.Cells.Item(100, 34).Value = 11231088300.5656
$ws.Cells.Item(100, 35).Value = 20056209198.9128
$ws.Cells.Item(100, 36).Value = 35815899284.6931
$ws.Cells.Item(100, 37).Value = 63959177372.4028
$ws.Cells.Item(100, 38).Value = 114216771094.808
$ws.Cells.Item(100, 39).Value = 203965581410.88
$ws.Cells.Item(100, 40).Value = 364236862953.741
$ws.Cells.Item(100, 41).Value = 650445488972.608
$ws.Cells.Item(100, 42).Value = 1161550016365.44
$ws.Cells.Item(100, 43).Value = 2074268272118.62
$ws.Cells.Item(100, 44).Value = 3704178730228.94
$ws.Cells.Item(100, 45).Value = 6614833891021.33

# Row 101
$ws.Cells.Item(101, 10).Value = 5534339.79452095
$ws.Cells.Item(101, 11).Value = 5239407.35483547
$ws.Cells.Item(101, 12).Value = 5305976.42178821
$ws.Cells.Item(101, 13).Value = 6067662.89541135
$ws.Cells.Item(101, 14).Value = 6355391.96059473
$ws.Cells.Item(101, 15).Value = 6318025.43231926
$ws.Cells.Item(101, 16).Value = 6318025
$ws.Cells.Item(101, 17).Value = 6489761.70275998
$ws.Cells.Item(101, 18).Value = 6666166.55657585
$ws.Cells.Item(101, 19).Value = 6847366.45123838
$ws.Cells.Item(101, 20).Value = 7033491.72565958
$ws.Cells.Item(101, 21).Value = 7224676.26162681
$ws.Cells.Item(101, 22).Value = 7421057.58010531
$ws.Cells.Item(101, 23).Value = 7622776.94015837
$ws.Cells.Item(101, 24).Value = 7829979.44055646
$ws.Cells.Item(101, 25).Value = 8042814.12414819
$ws.Cells.Item(101, 26).Value = 8261434.08506835
$ws.Cells.Item(101, 27).Value = 8485996.57886009
$ws.Cells.Item(101, 28).Value = 8716663.13559039
$ws.Cells.Item(101, 29).Value = 8953599.67604026
$ws.Cells.Item(101, 30).Value = 9196976.63105326
$ws.Cells.Item(101, 31).Value = 9446969.0641281
$ws.Cells.Item(101, 32).Value = 9703756.79734361
$ws.Cells.Item(101, 33).Value = 9967524.54070654
$ws.Cells.Item(101, 34).Value = 10238462.0250153
$ws.Cells.Item(101, 35).Value = 10516764.1383355
$ws.Cells.Item(101, 36).Value = 10802631.0661843
$ws.Cells.Item(101, 37).Value = 11096268.4355266
$ws.Cells.Item(101, 38).Value = 11397887.4626842
$ws.Cells.Item(101, 39).Value = 11707705.1052656
$ws.Cells.Item(101, 40).Value = 12025944.2182265
$ws.Cells.Item(101, 41).Value = 12352833.7141708
$ws.Cells.Item(101, 42).Value = 12688608.7280104
$ws.Cells.Item(101, 43).Value = 13033510.7860998
$ws.Cells.Item(101, 44).Value = 13387787.9799684
$ws.Cells.Item(101, 45).Value = 13751695.1447753

# Row 102
$ws.Cells.Item(102, 10).Value = 536578514.980824
$ws.Cells.Item(102, 11).Value = 597286241.842937
$ws.Cells.Item(102, 12).Value = 637099100.553869
$ws.Cells.Item(102, 13).Value = 656367407.571989
$ws.Cells.Item(102, 14).Value = 574373777.445775
$ws.Cells.Item(102, 15).Value = 457781451.868233
$ws.Cells.Item(102, 16).Value = 457781500
$ws.Cells.Item(102, 17).Value = 478465790.66445
$ws.Cells.Item(102, 18).Value = 500084675.4099
$ws.Cells.Item(102, 19).Value = 522680382.713528
$ws.Cells.Item(102, 20).Value = 546297049.094001
$ws.Cells.Item(102, 21).Value = 570980805.323975
$ws.Cells.Item(102, 22).Value = 596779866.537988
$ws.Cells.Item(102, 23).Value = 623744626.411778
$ws.Cells.Item(102, 24).Value = 651927755.596968
$ws.Cells.Item(102, 25).Value = 681384304.603405
$ws.Cells.Item(102, 26).Value = 712171811.33011
$ws.Cells.Item(102, 27).Value = 744350413.454879
$ws.Cells.Item(102, 28).Value = 777982965.902071
$ws.Cells.Item(102, 29).Value = 813135163.618033
$ws.Cells.Item(102, 30).Value = 849875669.893975
$ws.Cells.Item(102, 31).Value = 888276250.486968
$ws.Cells.Item(102, 32).Value = 928411913.801017
$ws.Cells.Item(102, 33).Value = 970361057.402055
$ws.Cells.Item(102, 34).Value = 1014205621.15303
$ws.Cells.Item(102, 35).Value = 1060031247.26822
$ws.Cells.Item(102, 36).Value = 1107927447.59937
$ws.Cells.Item(102, 37).Value = 1157987778.48051
$ws.Cells.Item(102, 38).Value = 1210310023.47288
$ws.Cells.Item(102, 39).Value = 1264996384.36691
$ws.Cells.Item(102, 40).Value = 1322153680.8145
$ws.Cells.Item(102, 41).Value = 1381893558.98135
$ws.Cells.Item(102, 42).Value = 1444332709.62703
$ws.Cells.Item(102, 43).Value = 1509593096.03866
$ws.Cells.Item(102, 44).Value = 1577802192.26363
$ws.Cells.Item(102, 45).Value = 1649093232.10641

# Row 103
$ws.Cells.Item(103, 10).Value = 1299102.98512899
$ws.Cells.Item(103, 11).Value = 1228137.20338854
$ws.Cells.Item(103, 12).Value = 1203104.99211033
$ws.Cells.Item(103, 13).Value = 1329584.91825699
$ws.Cells.Item(103, 14).Value = 1194752.03597448
$ws.Cells.Item(103, 15).Value = 941254.437436514
$ws.Cells.Item(103, 16).Value = 941254.4
$ws.Cells.Item(103, 17).Value = 940420.10532612
$ws.Cells.Item(103, 18).Value = 939586.550141589
$ws.Cells.Item(103, 19).Value = 938753.733790949
$ws.Cells.Item(103, 20).Value = 937921.655619323
$ws.Cells.Item(103, 21).Value = 937090.314972416
$ws.Cells.Item(103, 22).Value = 936259.711196513
$ws.Cells.Item(103, 23).Value = 935429.843638476
$ws.Cells.Item(103, 24).Value = 934600.711645749
$ws.Cells.Item(103, 25).Value = 933772.314566351
$ws.Cells.Item(103, 26).Value = 932944.651748882
$ws.Cells.Item(103, 27).Value = 932117.722542518
$ws.Cells.Item(103, 28).Value = 931291.526297012
$ws.Cells.Item(103, 29).Value = 930466.062362693
$ws.Cells.Item(103, 30).Value = 929641.330090467
$ws.Cells.Item(103, 31).Value = 928817.328831813
$ws.Cells.Item(103, 32).Value = 927994.057938789
$ws.Cells.Item(103, 33).Value = 927171.516764022
$ws.Cells.Item(103, 34).Value = 926349.704660717
$ws.Cells.Item(103, 35).Value = 925528.62098265
$ws.Cells.Item(103, 36).Value = 924708.265084171
$ws.Cells.Item(103, 37).Value = 923888.636320202
$ws.Cells.Item(103, 38).Value = 923069.734046237
$ws.Cells.Item(103, 39).Value = 922251.55761834
$ws.Cells.Item(103, 40).Value = 921434.106393147
$ws.Cells.Item(103, 41).Value = 920617.379727864
$ws.Cells.Item(103, 42).Value = 919801.376980267
$ws.Cells.Item(103, 43).Value = 918986.0975087
$ws.Cells.Item(103, 44).Value = 918171.540672078
$ws.Cells.Item(103, 45).Value = 917357.705829881

# Row 104
$ws.Cells.Item(104, 10).Value = 847891.146667666
$ws.Cells.Item(104, 11).Value = 802315.048000644
$ws.Cells.Item(104, 12).Value = 823092.738137487
$ws.Cells.Item(104, 13).Value = 938431.647770099
$ws.Cells.Item(104, 14).Value = 908454.125640036
$ws.Cells.Item(104, 15).Value = 868497.327264841
$ws.Cells.Item(104, 16).Value = 868497.3
$ws.Cells.Item(104, 17).Value = 898615.720978362
$ws.Cells.Item(104, 18).Value = 929778.611849987
$ws.Cells.Item(104, 19).Value = 962022.193549521
$ws.Cells.Item(104, 20).Value = 995383.943109194
$ws.Cells.Item(104, 21).Value = 1029902.63721874
$ws.Cells.Item(104, 22).Value = 1065618.39729592
$ws.Cells.Item(104, 23).Value = 1102572.73612005
$ws.Cells.Item(104, 24).Value = 1140808.60608272
$ws.Cells.Item(104, 25).Value = 1180370.44911175
$ws.Cells.Item(104, 26).Value = 1221304.24832652
$ws.Cells.Item(104, 27).Value = 1263657.58148456
$ws.Cells.Item(104, 28).Value = 1307479.67628169
$ws.Cells.Item(104, 29).Value = 1352821.46756982
$ws.Cells.Item(104, 30).Value = 1399735.65655905
$ws.Cells.Item(104, 31).Value = 1448276.77207279
$ws.Cells.Item(104, 32).Value = 1498501.23392715
$ws.Cells.Item(104, 33).Value = 1550467.41850828
$ws.Cells.Item(104, 34).Value = 1604235.72662377
$ws.Cells.Item(104, 35).Value = 1659868.65370711
$ws.Cells.Item(104, 36).Value = 1717430.8624568
$ws.Cells.Item(104, 37).Value = 1776989.25799436
$ws.Cells.Item(104, 38).Value = 1838613.06562888
$ws.Cells.Item(104, 39).Value = 1902373.91131825
$ws.Cells.Item(104, 40).Value = 1968345.90492069
$ws.Cells.Item(104, 41).Value = 2036605.72633341
$ws.Cells.Item(104, 42).Value = 2107232.71461841
$ws.Cells.Item(104, 43).Value = 2180308.96021901
$ws.Cells.Item(104, 44).Value = 2255919.4003744
$ws.Cells.Item(104, 45).Value = 2334151.91784306

# Row 111
$ws.Cells.Item(111, 10).Value = 33183.3692108051
$ws.Cells.Item(111, 11).Value = 29797.4021072387
$ws.Cells.Item(111, 12).Value = 30062.5076342299
$ws.Cells.Item(111, 13).Value = 32081.5924984002
$ws.Cells.Item(111, 14).Value = 34052.7325171703
$ws.Cells.Item(111, 15).Value = 33362.608580124
$ws.Cells.Item(111, 16).Value = 33362.61
$ws.Cells.Item(111, 17).Value = 35179.1927178324
$ws.Cells.Item(111, 18).Value = 37094.6877441062
$ws.Cells.Item(111, 19).Value = 39114.4808202274
$ws.Cells.Item(111, 20).Value = 41244.2509393822
$ws.Cells.Item(111, 21).Value = 43489.9863139954
$ws.Cells.Item(111, 22).Value = 45858.0012126131
$ws.Cells.Item(111, 23).Value = 48354.9537135467
$ws.Cells.Item(111, 24).Value = 50987.8644251971
$ws.Cells.Item(111, 25).Value = 53764.1362256946
$ws.Cells.Item(111, 26).Value = 56691.5750773548
$ws.Cells.Item(111, 27).Value = 59778.4119744747
$ws.Cells.Item(111, 28).Value = 63033.3260861793
$ws.Cells.Item(111, 29).Value = 66465.4691593876
$ws.Cells.Item(111, 30).Value = 70084.491250512
$ws.Cells.Item(111, 31).Value = 73900.5678582401
$ws.Cells.Item(111, 32).Value = 77924.4285336873
$ws.Cells.Item(111, 33).Value = 82167.3870483618
$ws.Cells.Item(111, 34).Value = 86641.3732047657
$ws.Cells.Item(111, 35).Value = 91358.9663790721
$ws.Cells.Item(111, 36).Value = 96333.4308901897
$ws.Cells.Item(111, 37).Value = 101578.753294661
$ws.Cells.Item(111, 38).Value = 107109.681712253
$ws.Cells.Item(111, 39).Value = 112941.767292818
$ws.Cells.Item(111, 40).Value = 119091.407941004
$ws.Cells.Item(111, 41).Value = 125575.894421767
$ws.Cells.Item(111, 42).Value = 132413.458976306
$ws.Cells.Item(111, 43).Value = 139623.326585128
$ws.Cells.Item(111, 44).Value = 147225.76902236
$ws.Cells.Item(111, 45).Value = 155242.161853305

# Row 112
$ws.Cells.Item(112, 10).Value = 191616.390889597
$ws.Cells.Item(112, 11).Value = 183780.936345654
$ws.Cells.Item(112, 12).Value = 187175.265813673
$ws.Cells.Item(112, 13).Value = 203489.457126293
$ws.Cells.Item(112, 14).Value = 198634.511450512
$ws.Cells.Item(112, 15).Value = 189245.572851842
$ws.Cells.Item(112, 16).Value = 189245.6
$ws.Cells.Item(112, 17).Value = 194245.395873519
$ws.Cells.Item(112, 18).Value = 199377.284428596
$ws.Cells.Item(112, 19).Value = 204644.755503007
$ws.Cells.Item(112, 20).Value = 210051.391134701
$ws.Cells.Item(112, 21).Value = 215600.867997689
$ws.Cells.Item(112, 22).Value = 221296.959902294
$ws.Cells.Item(112, 23).Value = 227143.540361453
$ws.Cells.Item(112, 24).Value = 233144.585224826
$ws.Cells.Item(112, 25).Value = 239304.175382486
$ws.Cells.Item(112, 26).Value = 245626.499540053
$ws.Cells.Item(112, 27).Value = 252115.857067136
$ws.Cells.Item(112, 28).Value = 258776.660921033
$ws.Cells.Item(112, 29).Value = 265613.440647675
$ws.Cells.Item(112, 30).Value = 272630.845461851
$ws.Cells.Item(112, 31).Value = 279833.647408814
$ws.Cells.Item(112, 32).Value = 287226.744609417
$ws.Cells.Item(112, 33).Value = 294815.164590979
$ws.Cells.Item(112, 34).Value = 302604.067706153
$ws.Cells.Item(112, 35).Value = 310598.750642123
$ws.Cells.Item(112, 36).Value = 318804.650022509
$ws.Cells.Item(112, 37).Value = 327227.346104433
$ws.Cells.Item(112, 38).Value = 335872.566573261
$ws.Cells.Item(112, 39).Value = 344746.190437601
$ws.Cells.Item(112, 40).Value = 353854.252027204
$ws.Cells.Item(112, 41).Value = 363202.945096491
$ws.Cells.Item(112, 42).Value = 372798.627036486
$ws.Cells.Item(112, 43).Value = 382647.823198039
$ws.Cells.Item(112, 44).Value = 392757.231329254
$ws.Cells.Item(112, 45).Value = 403133.726130163
